$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2025-05-24 Saturday" "2025-05-25 Sunday"
Replace-Text "65×38=2470" "85×89=7565"
Replace-Text "64×28=1792" "36×71=2556"
Replace-Text "60×92=5520" "52×37=1924"
Replace-Text "75×52=3900" "30×19=570"
Replace-Text "13×21=273" "15×44=660"
Replace-Text "93×55=5115" "48×34=1632"
Replace-Text "40×56=2240" "16×51=816"
Replace-Text "55×87=4785" "52×68=3536"
Replace-Text "39×87=3393" "61×70=4270"
Replace-Text "39×66=2574" "77×21=1617"
Replace-Text "84×14=1176" "46×33=1518"
Replace-Text "24×94=2256" "31×73=2263"
Replace-Text "44×50=2200" "43×99=4257"
Replace-Text "42×86=3612" "81×94=7614"
Replace-Text "34×70=2380" "16×61=976"
Replace-Text "52×97=5044" "29×78=2262"
Replace-Text "38×95=3610" "28×21=588"
Replace-Text "39×78=3042" "45×51=2295"
Replace-Text "33×52=1716" "99×88=8712"
Replace-Text "71×13=923" "86×50=4300"
Replace-Text "38×71=2698" "63×33=2079"
Replace-Text "41×54=2214" "81×23=1863"
Replace-Text "74×55=4070" "33×58=1914"
Replace-Text "75×92=6900" "31×83=2573"
Replace-Text "16×33=528" "11×96=1056"
